$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values per repull/recalculation of the data
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -4
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = -3
$ws.Range("F11").Value = 2
